$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54; everything from old row 54 downward
# shifts down by one (old row 54 -> new row 55, ..., old row 128 -> new row 129).
$ws.Rows("54").Insert()

# Populate the newly inserted row 54 with its data (same master/lookup
# columns as the surrounding Pina/Macroferia Regional de Talca rows, plus
# the specific observation values for this record).
$ws.Range("A54").Value = 5
$ws.Range("B54").Value = "Macroferia Regional de Talca"
$ws.Range("C54").Value = "Maule"
$ws.Range("D54").Value = 44413
$ws.Range("E54").Value = 7
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100108
$ws.Range("H54").Value = "Tropicales y subtropicales"
$ws.Range("I54").Value = 100108005
$ws.Range("J54").Value = "Piña"
$ws.Range("K54").Value = "Caramelo"
$ws.Range("L54").Value = "Segunda"
$ws.Range("M54").Value = 450
$ws.Range("N54").Value = 18000
$ws.Range("O54").Value = 18000
$ws.Range("P54").Value = 18000
$ws.Range("Q54").Value = "$/caja 14 unidades"
$ws.Range("R54").Value = "Ecuador"
$ws.Range("S54").Value = 1286
$ws.Range("T54").Value = 14
